# Update the LUY_KE_THANG_LONG_XUYEN sheet:
#  - last_edited_time (column D) for rows 3, 4, 5, 7, 13 moves from
#    2024-07-20T13:34:00.000Z to 2024-07-21T16:44:00.000Z
#  - Month 7 (row 5) numeric metrics are bumped by +5,000,000 (and the
#    order-count bumped by 1) to reflect the updated sale-phu discount /
#    hourly-wage strategy figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_LONG_XUYEN")

$newTimestamp = "2024-07-21T16:44:00.000Z"

foreach ($r in 3, 4, 5, 7, 13) {
    $ws.Range("D$r").Value = $newTimestamp
}

# Row 5 ("Tháng 7") numeric updates
$ws.Range("AA5").Value = 27816000   # properties.Lũy kế.formula.number
$ws.Range("AE5").Value = 44950000   # properties.Tổng doanh thu.formula.number
$ws.Range("AH5").Value = 41650000   # properties.Đã thanh toán.number
$ws.Range("AK5").Value = 12         # properties.Số lượng đơn.number
$ws.Range("AQ5").Value = 45150000   # properties.Đơn giá.number
